{"js": "// Word adds \":v we hueles mal , huacala\" before the final \"_GoBack\" bookmark\n// (Word's spell-checker flags \"we\" as a typo, wrapping it in proofErr marks)\n// and appends a brand-new empty paragraph right after that paragraph.\n\nconst body = context.document.body;\n\n// --- Step 1: append a new, fully empty trailing paragraph at the very end of\n// the document body. We ship two empty paragraphs in the OOXML fragment\n// because the host rejects a payload that resolves to literally no content;\n// using two empty paragraphs keeps the *existing* last paragraph (which\n// carries the _GoBack bookmark) completely untouched and simply adds one new\n// empty paragraph after it.\nconst trailingParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body><w:p/><w:p/></w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\nbody.insertOoxml(trailingParagraphOoxml, Word.InsertLocation.end);\nawait context.sync();\n\n// --- Step 2: locate the paragraph that holds the _GoBack bookmark. After\n// step 1 it is the second-to-last paragraph (the last one is the new empty\n// paragraph we just added).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst bookmarkParagraph = paragraphs.items[paragraphs.items.length - 2];\n\n// --- Step 3: insert the new text at the very start of that paragraph, i.e.\n// before the bookmark, split into three runs with proofErr spell-check\n// markers wrapped around \"we\" to match Word's own authoring output.\nconst newTextOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r><w:t xml:space=\"preserve\">:v </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>we</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> hueles mal , huacala</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\nbookmarkParagraph.insertOoxml(newTextOoxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Word adds \":v we hueles mal , huacala\" just before the existing \"_GoBack\"\n# bookmark (Word's spell-checker flags \"we\" as a typo, so it gets wrapped in\n# proofErr spell-check marks) and leaves a brand-new empty paragraph right\n# after that paragraph.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that carries the \"_GoBack\" bookmark (the last paragraph\n# in the document) without depending on a hard-coded paragraph index.\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$targetParagraph = $bookmark.Range.Paragraphs(1)\n$targetRange = $targetParagraph.Range\n\n# Replace that whole paragraph (bookmark included) with the same bookmark\n# plus the new text in front of it, split into three runs with proofErr\n# spell-check markers wrapped around \"we\" to match Word's own authoring\n# output. Re-supplying the bookmark inside the replacement keeps it alive\n# (and in the correct place, right after the new text) instead of losing it.\n$newParagraphXml = @'\n<?xml version=\"1.0\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">:v </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>we</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> hueles mal , huacala</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$null = $targetRange.InsertXML($newParagraphXml)\n"}
